# Apply updated dSF ("F" column) values to Sheet1.
# Changes: F6 -1->-4, F12 4->1, F13 -1->-2, F15 -9->-10, F20 3->2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F6").Value = -4
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -2
$ws.Range("F15").Value = -10
$ws.Range("F20").Value = 2
